$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row, SpdLin, AngSpd, sp_left, sp_right, Filename)
# sp_left/sp_right/Filename are text values (some look numeric, e.g. "0",
# "160", "200") so they are entered with a leading apostrophe to force
# Excel to store them as text (shared strings) instead of numbers.
$data = @(
    @(13, -1583.58,  -98.08, "'0",   "'0",   "'26_05_2017_1302"),
    @(14, 24606.82,  -41.52, "'160", "'160", "'26_05_2017_1304"),
    @(15, 14071.75,    7.63, "'160", "'160", "'26_05_2017_1318"),
    @(16,  3954.1,   -23.1,  "'160", "'160", "'26_05_2017_1319"),
    @(17, 12999.59,  -83.64, "'160", "'160", "'26_05_2017_1321"),
    @(18,  -380.76,   15.63, "'160", "'200", "'26_05_2017_1322"),
    @(19, 23743.07,  -69.39, "'160", "'200", "'26_05_2017_1323"),
    @(20, 28092.45,  -21.95, "'160", "'200", "'26_05_2017_1324")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# Match the formatting already used by the existing data rows (row 3..12)
# by copying row 12's cell formats down over the newly added rows.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
